# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-01-25 Saturday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2025-01-26 Sunday", 2)

# Update the multiplication problems in the table. Each problem cell is
# addressed directly by (row, column) so that identical "old" / "new"
# strings appearing elsewhere in the table (e.g. 452x9= becomes 559x7=,
# while a different cell's original 559x7= becomes 887x4=) can never be
# confused with a global find/replace.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "265×6="   # was 554×2=
$t.Cell(1,2).Range.Text  = "426×7="   # was 467×3=
$t.Cell(1,3).Range.Text  = "915×5="   # was 683×3=
$t.Cell(1,4).Range.Text  = "107×2="   # was 716×4=
$t.Cell(1,5).Range.Text  = "273×2="   # was 835×6=

$t.Cell(5,1).Range.Text  = "173×9="   # was 764×2=
$t.Cell(5,2).Range.Text  = "605×5="   # was 578×7=
$t.Cell(5,3).Range.Text  = "274×4="   # was 890×4=
$t.Cell(5,4).Range.Text  = "708×2="   # was 775×3=
$t.Cell(5,5).Range.Text  = "426×2="   # was 872×6=

$t.Cell(10,1).Range.Text = "232×6="   # was 583×6=
$t.Cell(10,2).Range.Text = "526×4="   # was 141×3=
$t.Cell(10,3).Range.Text = "548×5="   # was 337×5=
$t.Cell(10,4).Range.Text = "778×2="   # was 473×5=
$t.Cell(10,5).Range.Text = "559×7="   # was 452×9=

$t.Cell(15,1).Range.Text = "803×7="   # was 478×8=
$t.Cell(15,2).Range.Text = "851×4="   # was 473×8=
$t.Cell(15,3).Range.Text = "123×6="   # was 113×6=
$t.Cell(15,4).Range.Text = "457×7="   # was 475×9=
$t.Cell(15,5).Range.Text = "419×3="   # was 950×3=

$t.Cell(20,1).Range.Text = "809×9="   # was 593×3=
$t.Cell(20,2).Range.Text = "887×4="   # was 559×7=
$t.Cell(20,3).Range.Text = "895×4="   # was 131×3=
$t.Cell(20,4).Range.Text = "961×9="   # was 916×4=
$t.Cell(20,5).Range.Text = "675×6="   # was 506×8=
